# "add stats vs Oplati" - update standings table with the latest round's
# results and append the new round's match rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Standings table (rows 5-18): team order re-ranked + updated
#    games / wins / losses / points. Column C (team) and G (score)
#    already hold the right values for the team that ends up in each
#    row, so only the team name, score text and W/L/Pts numbers need
#    to be (re)written per row.
# ---------------------------------------------------------------------

$standings = @(
    # place, team,                      games, wins, losses, score,        points
    @(1,  "GOLDEN HILL",                 5, 5, 0, "357 - 306", 10),
    @(2,  "ISsoft",                      5, 4, 1, "291 - 256", 9),
    @(3,  "Эра-Недвижимости плюс",       5, 4, 1, "373 - 310", 9),
    @(4,  "Грушвиль",                    5, 4, 1, "391 - 332", 9),
    @(5,  "ОПЛАТИ",                      5, 3, 2, "318 - 288", 8),
    @(6,  "БГУФК",                       5, 2, 3, "309 - 273", 7),
    @(7,  "Mapogo males",                4, 3, 1, "296 - 284", 7),
    @(8,  "VSS",                         5, 2, 3, "310 - 337", 7),
    @(9,  "SIRIUS",                      5, 2, 3, "306 - 294", 7),
    @(10, "NORD",                        5, 2, 3, "254 - 395", 7),
    @(11, "Стрела",                      5, 1, 4, "296 - 323", 6),
    @(12, "Eagles",                      5, 1, 4, "280 - 294", 6),
    @(13, "ЛФК",                         5, 1, 4, "293 - 337", 6),
    @(14, "Минск 7х",                    4, 0, 4, "190 - 235", 4)
)

$row = 5
foreach ($team in $standings) {
    $ws.Range("C$row").Value = $team[1]
    $ws.Range("D$row").Value = $team[2]
    $ws.Range("E$row").Value = $team[3]
    $ws.Range("F$row").Value = $team[4]
    $ws.Range("G$row").Value = $team[5]
    $ws.Range("H$row").Value = $team[6]
    $row++
}

# ---------------------------------------------------------------------
# 2. Append the new round (two match days) at the bottom of the sheet.
#    Re-use the formatting of the previous round by copying a date row
#    and a match row, then overwrite their contents.
# ---------------------------------------------------------------------

$newRows = @(
    @(45626, $null),
    @($null, "Грушвиль - БГУФК 74:70 (16:30, БНТУ)"),
    @($null, "Стрела - ISsoft 42:56 (18:00, БНТУ)"),
    @($null, "Mapogo males - ЛФК 69:53 (19:30, БНТУ)"),
    @(45627, $null),
    @($null, "Минск 7х - NORD 41:53 (11:00, БНТУ)"),
    @($null, "Эра-Недвижимости плюс - VSS 80:64 (12:30, БНТУ)"),
    @($null, "GOLDEN HILL - Eagles 72:66 (14:00, БНТУ)"),
    @($null, "SIRIUS - ОПЛАТИ 49:60 (15:30, БНТУ)")
)

$destRow = 55
foreach ($entry in $newRows) {
    $date = $entry[0]
    $text = $entry[1]

    if ($null -ne $date) {
        $ws.Range("B50:H50").Copy()
        $ws.Range("B${destRow}:H${destRow}").PasteSpecial()
        $ws.Range("B$destRow").Value = $date
    } else {
        $ws.Range("B49:H49").Copy()
        $ws.Range("B${destRow}:H${destRow}").PasteSpecial()
        $ws.Range("B$destRow").Value = $text
    }

    $destRow++
}

$excel.CutCopyMode = $false
